$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 26 previously duplicated "Cebus imitator" (the real species for this
# accession is Cebus capucinus). Fix the species name; this introduces a new
# shared string.
$ws.Range("A26").Value = "Cebus capucinus"

# Clear the (redundant) explicit cell style that had accumulated on B8:B33
# so they fall back to the default/normal style.
$ws.Range("B8:B33").ClearFormats()

$ws.Range("F30").Select()
